$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly values between row 2 and row 3 for columns D (Fecha),
# M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio
# ponderado) and S (Precio $/Kg).

$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $addr2 = "$col`2"
    $addr3 = "$col`3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}
